$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (everything below shifts down by one row)
$ws.Rows("2:2").Insert()

# New journal entry content
$ws.Range("A2").Value = "mercredi 2 février 2016"
$ws.Range("B2").Value = "Création procédure stockée`nGUI`nTests sur la bd`nAjout création de compte`nAjout checkbox dynamiques GUI"

# Row height for the new row
$ws.Range("A2:B2").RowHeight = 92.25

# Borders: thin all around, then thicken the top edge to medium
$rng = $ws.Range("A2:B2")
$rng.Borders.LineStyle = 1
$rng.Borders.Item(8).Weight = -4138

# Alignment/wrap to match the rest of the table
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4108
$ws.Range("B2").HorizontalAlignment = -4131
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("B2").WrapText = $true

# Move the selection as in the authored workbook
$ws.Range("C14").Select()
